$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.73180333333334
$ws.Range("H2").Value = 68.19541000000001
$ws.Range("I2").Value = 0.007290521456144306
$ws.Range("J2").Value = 0.007290521456144306
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 42.1577583942389
$ws.Range("R2").Value = 379.4198255481501
$ws.Range("S2").Value = 0.0001203936041835694
$ws.Range("T2").Value = 0.0001203936041835694

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.73180333333334
$ws.Range("H3").Value = 68.19541000000001
$ws.Range("I3").Value = 0.007290521456144306
$ws.Range("J3").Value = 0.007290521456144306
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 2112.096508512881
$ws.Range("R3").Value = 19008.86857661593
$ws.Range("S3").Value = 0.006031699044941344
$ws.Range("T3").Value = 0.006031699044941345

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 22.73180333333334
$ws.Range("H4").Value = 68.19541000000001
$ws.Range("I4").Value = 0.007290521456144306
$ws.Range("J4").Value = 0.007290521456144306
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 3.685643665253334
$ws.Range("R4").Value = 33.17079298728001
$ws.Range("S4").Value = 0.00001052541552249201
$ws.Range("T4").Value = 0.00001052541552249201

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 22.73180333333334
$ws.Range("H5").Value = 68.19541000000001
$ws.Range("I5").Value = 0.007290521456144306
$ws.Range("J5").Value = 0.007290521456144306
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 394.9535275833056
$ws.Range("R5").Value = 3554.58174824975
$ws.Range("S5").Value = 0.0011279033914969
$ws.Range("T5").Value = 0.0011279033914969

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3042.696044666667
$ws.Range("H6").Value = 9128.088134
$ws.Range("I6").Value = 0.9758504625824999
$ws.Range("J6").Value = 0.9758504625824997
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 5642.897874717534
$ws.Range("R6").Value = 50786.08087245781
$ws.Range("S6").Value = 0.01611491784795388
$ws.Range("T6").Value = 0.01611491784795388

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3042.696044666667
$ws.Range("H7").Value = 9128.088134
$ws.Range("I7").Value = 0.9758504625824999
$ws.Range("J7").Value = 0.9758504625824997
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 282708.2215242824
$ws.Range("R7").Value = 2544373.993718541
$ws.Range("S7").Value = 0.8073546369174731
$ws.Range("T7").Value = 0.8073546369174731

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3042.696044666667
$ws.Range("H8").Value = 9128.088134
$ws.Range("I8").Value = 0.9758504625824999
$ws.Range("J8").Value = 0.9758504625824997
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 493.3305658980748
$ws.Range("R8").Value = 4439.975093082672
$ws.Range("S8").Value = 0.001408847318848566
$ws.Range("T8").Value = 0.001408847318848566

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3042.696044666667
$ws.Range("H9").Value = 9128.088134
$ws.Range("I9").Value = 0.9758504625824999
$ws.Range("J9").Value = 0.9758504625824997
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 52865.29707226063
$ws.Range("R9").Value = 475787.6736503456
$ws.Range("S9").Value = 0.1509720604982242
$ws.Range("T9").Value = 0.1509720604982242

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.93578
$ws.Range("H10").Value = 5.80734
$ws.Range("I10").Value = 0.0006208414447999517
$ws.Range("J10").Value = 0.0006208414447999516
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 3.5900427409
$ws.Range("R10").Value = 32.3103846681
$ws.Range("S10").Value = 0.00001025239958700167
$ws.Range("T10").Value = 0.00001025239958700167

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.93578
$ws.Range("H11").Value = 5.80734
$ws.Range("I11").Value = 0.0006208414447999517
$ws.Range("J11").Value = 0.0006208414447999516
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 179.86052928998
$ws.Range("R11").Value = 1618.74476360982
$ws.Range("S11").Value = 0.0005136434714836331
$ws.Range("T11").Value = 0.000513643471483633

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.93578
$ws.Range("H12").Value = 5.80734
$ws.Range("I12").Value = 0.0006208414447999517
$ws.Range("J12").Value = 0.0006208414447999516
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 0.3138596260800001
$ws.Range("R12").Value = 2.82473663472
$ws.Range("S12").Value = 0.0000008963164321526731
$ws.Range("T12").Value = 0.0000008963164321526731

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.93578
$ws.Range("H13").Value = 5.80734
$ws.Range("I13").Value = 0.0006208414447999517
$ws.Range("J13").Value = 0.0006208414447999516
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 33.6331934785
$ws.Range("R13").Value = 302.6987413065
$ws.Range("S13").Value = 0.00009604925729716422
$ws.Range("T13").Value = 0.00009604925729716419

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 50.63053333333333
$ws.Range("H14").Value = 151.8916
$ws.Range("I14").Value = 0.016238174516556
$ws.Range("J14").Value = 0.016238174516556
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 93.89795258822221
$ws.Range("R14").Value = 845.0815732939999
$ws.Range("S14").Value = 0.0002681526098194738
$ws.Range("T14").Value = 0.0002681526098194739

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 50.63053333333333
$ws.Range("H15").Value = 151.8916
$ws.Range("I15").Value = 0.016238174516556
$ws.Range("J15").Value = 0.016238174516556
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 4704.271416982977
$ws.Range("R15").Value = 42338.44275284679
$ws.Range("S15").Value = 0.01343440003740153
$ws.Range("T15").Value = 0.01343440003740153

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 50.63053333333333
$ws.Range("H16").Value = 151.8916
$ws.Range("I16").Value = 0.016238174516556
$ws.Range("J16").Value = 0.016238174516556
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 8.209032152533334
$ws.Range("R16").Value = 73.8812893728
$ws.Range("S16").Value = 0.00002344325232997568
$ws.Range("T16").Value = 0.00002344325232997568

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 50.63053333333333
$ws.Range("H17").Value = 151.8916
$ws.Range("I17").Value = 0.016238174516556
$ws.Range("J17").Value = 0.016238174516556
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 879.6797794788888
$ws.Range("R17").Value = 7917.118015309999
$ws.Range("S17").Value = 0.002512178617005022
$ws.Range("T17").Value = 0.002512178617005022
